$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$result = $ws.Cells.Replace("na", "NA", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
Write-Host "Replace result: $result"
